# Update the cached text of the "date" and "slide number" placeholder
# fields on the slide master and every slide layout.
#
#   datetimeFigureOut field: 28.10.2023 -> 16.02.2025
#   slidenum field:          ‹№›        -> ‹#›

$p = $ppt.ActivePresentation

$newDate = "16.02.2025"
$newSlideNum = [string][char]0x2039 + [char]0x23 + [char]0x203A

function Update-HeadersFooters($headersFooters) {
    $headersFooters.DateAndTime.Text = $newDate
    $headersFooters.SlideNumber.Text = $newSlideNum
}

# Slide master
Update-HeadersFooters $p.SlideMaster.HeadersFooters

# Every custom (slide) layout attached to the master
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-HeadersFooters $layout.HeadersFooters
}
